$wb = $excel.ActiveWorkbook

# --- Add "Netherlands" sheet as a copy of "Greece", placed right after it ---
$greece = $wb.Worksheets.Item("Greece")
$greece.Copy($null, $greece)
$netherlands = $wb.Worksheets.Item($wb.Worksheets.Count)
$netherlands.Name = "Netherlands"
# Allocate shared strings in the same order as the target file: B4 before B2.
$netherlands.Range("B4").Value = "NGC-3144/T2175"
$netherlands.Range("B2").Value = "Netherlands Market"

# --- Add "Austria" sheet as a copy of "Greece", placed after "Netherlands" ---
$greece.Copy($null, $netherlands)
$austria = $wb.Worksheets.Item($wb.Worksheets.Count)
$austria.Name = "Austria"
# Row 11 ("MX Minerva Bridge Kit") is missing on this sheet - delete it, shifting rows up.
$austria.Rows.Item(11).Delete()
$austria.Range("B4").Value = "NGC-3817/T2271"
# NOTE: B2 keeps the copied "Netherlands Market" text on this sheet (matches source data).
$austria.Range("B2").Value = "Netherlands Market"

# --- Add "Denmark" sheet as a copy of "Greece", placed after "Austria" ---
$greece.Copy($null, $austria)
$denmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Name = "Denmark"
$denmark.Range("B4").Value = "NGC-2913/T2748"
$denmark.Range("B2").Value = "Denmark Market"

# --- Selections on each new sheet ---
$austria.Range("B4").Select()
$denmark.Range("B4").Select()
# Netherlands ends up as the active / tab-selected sheet.
$netherlands.Range("E22").Select()
